$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rename the header labels in row 1.
#    "<Name>_old" -> "<Name>_FV2310"
#    "<Name>_new" -> "<Name>_FV2404"
#    ("diff" in K1 stays untouched.)
# ---------------------------------------------------------------------------
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")
$oldLetters = @("A", "B", "C", "D", "E", "F", "G", "H", "I", "J")
$newLetters = @("L", "M", "N", "O", "P", "Q", "R", "S", "T", "U")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $ws.Range($oldLetters[$i] + "1").Value = $baseNames[$i] + "_FV2310"
    $ws.Range($newLetters[$i] + "1").Value = $baseNames[$i] + "_FV2404"
}

# ---------------------------------------------------------------------------
# 2) Turn the used range A1:U55 into an Excel Table ("Table1").
# ---------------------------------------------------------------------------
$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U55"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# ---------------------------------------------------------------------------
# 3) Freeze the header row (split below row 1).
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
